$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Minimize the workbook window (bookViews/workbookView minimized="1")
$win = $wb.Windows.Item(1)
$win.WindowState = -4140

# Clear the stray zero values that were left in the sheet (keep style/formatting)
$ws.Range("J4").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("J5").ClearContents()

# Move/replace the active selection from M1 to H8
$ws.Range("H8").Select()
